# Generate Report for Handoff
# -------------------------------------------------------------------
# This localization-status workbook tracked exactly one source file
# (992a7d4b-...md). A new handoff round adds two image dependencies
# (.png files) referenced by a brand-new markdown source file, and the
# original row's data is replaced by the first of those new entries.
#
#   Sheet "Overview" : File Name | zh-cn | de-de | Latest Handoff Date
#   Sheet "zh-cn"     : Source File Name | File Extension | Status | ...
#   Sheet "de-de"     : Source File Name | File Extension | Status | ...
# -------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$repoCommit      = "adc862e90094ba9fe3dc41b06f0835493585d693"
$handoffZhCommit = "060f078d458a78e31d64f074b3828df8b93f8512"
$handoffDeCommit = "104291988fea2125ed4e0c3931e0a43636ff7c0a"

function SourceUrl([string]$fileName) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/$repoCommit/e2e/$fileName"
}
function HandoffZhUrl([string]$fileName) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffZhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$fileName"
}
function HandoffDeUrl([string]$fileName) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffDeCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$fileName"
}

# The three source files now being reported on.
$png1 = "0ad1acc5-38d9-45b1-9b10-09866b7c3bea.png"
$png2 = "cd1519f1-b8e0-4846-8b85-e3d8ba53707b.png"
$md   = "e3859f52-495a-4eb6-9607-0997697abeca.md"

$png1Target = "7b43e3b496085599f1cf59966f6a278b2b14e414.png"
$png2Target = "410d1853730b7d7fb7fde9013a08f1f9383a8da4.png"
$mdTargetZh = "e3859f52-495a-4eb6-9607-0997697abeca.61d2b984721d9b2806e022a78172838799ee1ae7.zh-cn.xlf"
$mdTargetDe = "e3859f52-495a-4eb6-9607-0997697abeca.61d2b984721d9b2806e022a78172838799ee1ae7.de-de.xlf"

$readyForHandoff = "Ready for handoff"
$latestHandoffDate = "2016-46-19 06:46:57"

# =====================================================================
# Sheet "Overview"
# =====================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Range("A2").Value = $png1
$wsOverview.Range("B2").Value = $readyForHandoff
$wsOverview.Range("C2").Value = $readyForHandoff
$wsOverview.Range("D2").Value = $latestHandoffDate

$wsOverview.Range("A3").Value = $png2
$wsOverview.Range("B3").Value = $readyForHandoff
$wsOverview.Range("C3").Value = $readyForHandoff
$wsOverview.Range("D3").Value = $latestHandoffDate

$wsOverview.Range("A4").Value = $md
$wsOverview.Range("B4").Value = $readyForHandoff
$wsOverview.Range("C4").Value = $readyForHandoff
$wsOverview.Range("D4").Value = $latestHandoffDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), (SourceUrl $png1), "", "", $png1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), (SourceUrl $png2), "", "", $png2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), (SourceUrl $md),   "", "", $md)

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("B2").Hyperlinks.Delete()
$wsZh.Range("D2").Hyperlinks.Delete()

# Row 2 - first png (dependency)
$wsZh.Range("A2").Value = $png1
$wsZh.Range("B2").Value = ".png"
$wsZh.Range("C2").Value = $readyForHandoff
$wsZh.Range("D2").Value = $png1Target
$wsZh.Range("E2").Value = "2016-03-19 06:46:54"
$wsZh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H2").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = "IsDependency"
$wsZh.Range("J2").Value = "e2e\" + $md

# Row 3 - second png (dependency)
$wsZh.Range("A3").Value = $png2
$wsZh.Range("B3").Value = ".png"
$wsZh.Range("C3").Value = $readyForHandoff
$wsZh.Range("D3").Value = $png2Target
$wsZh.Range("E3").Value = "2016-03-19 06:46:54"
$wsZh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "IsDependency"
$wsZh.Range("J3").Value = "e2e\" + $md

# Row 4 - the markdown source itself (included directly)
$wsZh.Range("A4").Value = $md
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $readyForHandoff
$wsZh.Range("D4").Value = $mdTargetZh
$wsZh.Range("E4").Value = "2016-03-19 06:46:54"
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), (SourceUrl $png1), "", "", $png1)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), (SourceUrl $png1), "", "", ".png")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), (HandoffZhUrl $png1Target), "", "", $png1Target)

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), (SourceUrl $png2), "", "", $png2)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), (SourceUrl $png2), "", "", ".png")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), (HandoffZhUrl $png2Target), "", "", $png2Target)

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), (SourceUrl $md), "", "", $md)
$wsZh.Hyperlinks.Add($wsZh.Range("B4"), (SourceUrl $md), "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D4"), (HandoffZhUrl $mdTargetZh), "", "", $mdTargetZh)

# =====================================================================
# Sheet "de-de"
# =====================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("B2").Hyperlinks.Delete()
$wsDe.Range("D2").Hyperlinks.Delete()

# Row 2 - first png (dependency)
$wsDe.Range("A2").Value = $png1
$wsDe.Range("B2").Value = ".png"
$wsDe.Range("C2").Value = $readyForHandoff
$wsDe.Range("D2").Value = $png1Target
$wsDe.Range("E2").Value = "2016-03-19 06:46:57"
$wsDe.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H2").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = "IsDependency"
$wsDe.Range("J2").Value = "e2e\" + $md

# Row 3 - second png (dependency)
$wsDe.Range("A3").Value = $png2
$wsDe.Range("B3").Value = ".png"
$wsDe.Range("C3").Value = $readyForHandoff
$wsDe.Range("D3").Value = $png2Target
$wsDe.Range("E3").Value = "2016-03-19 06:46:57"
$wsDe.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "IsDependency"
$wsDe.Range("J3").Value = "e2e\" + $md

# Row 4 - the markdown source itself (included directly)
$wsDe.Range("A4").Value = $md
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $readyForHandoff
$wsDe.Range("D4").Value = $mdTargetDe
$wsDe.Range("E4").Value = "2016-03-19 06:46:57"
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), (SourceUrl $png1), "", "", $png1)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), (SourceUrl $png1), "", "", ".png")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), (HandoffDeUrl $png1Target), "", "", $png1Target)

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), (SourceUrl $png2), "", "", $png2)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), (SourceUrl $png2), "", "", ".png")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), (HandoffDeUrl $png2Target), "", "", $png2Target)

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), (SourceUrl $md), "", "", $md)
$wsDe.Hyperlinks.Add($wsDe.Range("B4"), (SourceUrl $md), "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D4"), (HandoffDeUrl $mdTargetDe), "", "", $mdTargetDe)

Write-Host "Report generated for handoff."
